# Workbook: 832491-奥迪威.xlsx
# Add a new "2022-Q3" sheet (holding the Q3 fund-holding breakdown) positioned
# right after "总计" and before the existing "2022-Q2" sheet, and append a new
# "2022-Q2" summary row on the "总计" sheet (which previously only had one row).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right before the existing "2022-Q2"
#    sheet (i.e. right after "总计"), then populate it with the Q3 fund data.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# Header row - reuse the exact look of the "总计" header (bold + thin border,
# centered/top aligned) by copying its formatting over.
$wsTotal.Range("B1:D1").Copy() | Out-Null
$wsQ3.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$wsTotal.Range("A2").Copy() | Out-Null
$wsQ3.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$wsQ3.Application.CutCopyMode = 0

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Row 2
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").NumberFormat = "@"
$wsQ3.Range("B2").Value = "'014269"
$wsQ3.Range("C2").Value = "嘉实北交所精选两年定期混合A"
$wsQ3.Range("D2:F2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "'2.72"
$wsQ3.Range("E2").Value = "'90.37"
$wsQ3.Range("F2").Value = "'5.06"
$wsQ3.Range("G2").NumberFormat = "@"
$wsQ3.Range("G2").Value = "'0.1376"
$wsQ3.Range("H2").Value = 8

# Row 3
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").NumberFormat = "@"
$wsQ3.Range("B3").Value = "'014270"
$wsQ3.Range("C3").Value = "嘉实北交所精选两年定期混合C"
$wsQ3.Range("D3:F3").NumberFormat = "@"
$wsQ3.Range("D3").Value = "'0.53"
$wsQ3.Range("E3").Value = "'90.37"
$wsQ3.Range("F3").Value = "'5.06"
$wsQ3.Range("G3").NumberFormat = "@"
$wsQ3.Range("G3").Value = "'0.0268"
$wsQ3.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: the existing row 2 ("2022-Q2") now represents
#    the latest quarter ("2022-Q3" with its own totals), and a new row 3 is
#    appended preserving the original "2022-Q2" totals.
# ---------------------------------------------------------------------------
$wsTotal.Range("A2").Copy() | Out-Null
$wsTotal.Range("A3").PasteSpecial(-4122) | Out-Null
$wsTotal.Application.CutCopyMode = 0

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.22

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.16
